$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "input" to "expected_output"
$ws.Name = "expected_output"

# Clear the old formatted (date) cell E2 that is no longer used
$ws.Cells.Item(2, 5).Clear()

# Write the header row A1:F1 with the expected column names
$headers = @("tool_pid", "tool_code", "tool_type", "tool_price", "tool_inventory_start_date", "tool_inventory_end_date")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Update the selection to match the target (F2)
$ws.Range("F2").Select()
